$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Epoch")

# Insert 20 new rows right after the header row (row 1), shifting existing data down
$ws.Rows("2:21").Insert()

$newValues = @(
    "Epoch:521, time:9.454063, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:520, time:9.624205, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:519, time:10.734310, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:518, time:10.198558, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:517, time:9.551558, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:516, time:10.013921, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:515, time:9.508407, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:514, time:9.746659, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:513, time:9.784823, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:512, time:9.676936, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:511, time:9.462044, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:510, time:9.686640, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:509, time:9.553230, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:508, time:9.504339, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:507, time:9.519701, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:506, time:9.741438, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:505, time:10.239930, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:504, time:9.770286, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:503, time:9.520222, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53",
    "Epoch:502, time:9.511056, test_Acc: 38.30, test_bacc: 33.49, test_f1: 25.53"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
